$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147, shifting existing rows 147:159 down to 148:160.
$ws.Rows.Item(147).Insert()

# Fill the newly inserted row 147 with the new record's data.
$ws.Cells.Item(147, 1).Value = 1
$ws.Cells.Item(147, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(147, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(147, 4).Value = 44946
$ws.Cells.Item(147, 5).Value = 15
$ws.Cells.Item(147, 6).Value = "Fruta"
$ws.Cells.Item(147, 7).Value = 100106
$ws.Cells.Item(147, 8).Value = "Oleaginosos"
$ws.Cells.Item(147, 9).Value = 100106002
$ws.Cells.Item(147, 10).Value = "Palta"
$ws.Cells.Item(147, 11).Value = "Hass"
$ws.Cells.Item(147, 12).Value = "Segunda"
$ws.Cells.Item(147, 13).Value = 400
$ws.Cells.Item(147, 14).Value = 37000
$ws.Cells.Item(147, 15).Value = 38000
$ws.Cells.Item(147, 16).Value = 37500
$ws.Cells.Item(147, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(147, 18).Value = "Perú"
$ws.Cells.Item(147, 19).Value = 3750
$ws.Cells.Item(147, 20).Value = 10
